$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("K7").Value = "2016-08-27 04:53:51"
Write-Host "done"
